# Target worksheet: "SAB-BCV01-B03 02" (the active sheet / tabSelected sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new blank columns at I so the existing "Color Code" column
#    (I) shifts right to K, carrying its data, styles and the exact
#    "best fit" column width along with it. This creates fresh, empty
#    I and J columns (inheriting the row's existing cell style) ready to
#    receive the new "EW"/"LR" (motor position) data.
$ws.Range("I1:I1048576").Insert()
$ws.Range("I1:I1048576").Insert()

# 2) Header row (row 3): new column headers "EW" and "LR"
$ws.Range("I3").Value = "EW"
$ws.Range("J3").Value = "LR"

# 3) Data rows: column J gets the motor position (L/R), column I gets a
#    quantity only for the rows that need it (rows 5 and 8); otherwise it
#    stays blank (as left behind by the column insert above).
$ws.Range("J4").Value = "L"

$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = "R"

$ws.Range("J6").Value = "L"

$ws.Range("J7").Value = "R"

$ws.Range("I8").Value = 1500
$ws.Range("J8").Value = "L"

$ws.Range("J9").Value = "R"

# 4) Row 9's new J cell picked up the "fill" flavoured style (s=2) from the
#    row's other inserted cell; the target file keeps J9 on the plain
#    border-only style (s=1, same as the rest of the J column). Re-apply
#    that style by copying the format from a known s=1 cell.
$ws.Range("I4").Copy()
$ws.Range("J9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 5) Restore the cursor/selection to match the authored workbook.
$ws.Range("J14").Select()
